# Apply the "Abandonded Chapel" quest re-ordering + "Purgatory's Lantern" rename edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row at 49 - this shifts the existing rows 49-56 down to 50-57,
#    which moves the old "Alchemy Corrupted Church" row from 54 to 55 (so its
#    required_quest_item_id lookup below still targets D55), and moves the old
#    "Abandonded Chapel" row from 56 to 57 (it will be removed once its data is
#    copied into the freshly inserted row 49).
$ws.Rows.Item(49).Insert()

# 2) Populate the new row 49 with the "Abandonded Chapel" quest data (this is
#    the same data that used to live in row 56, plus a game_map_id that it
#    did not have before).
$ws.Range("A49").Value = "Abandonded Chapel"
$ws.Range("B49").Value = "Delusional Memories"
$ws.Range("E49").Value = "An old decrepid chapel in the middle of no where. Half burned, half rotted, what remains is a story of the past."
$ws.Range("G49").Value = 1
$ws.Range("H49").Value = 3
$ws.Range("I49").Value = 1
$ws.Range("J49").Value = 208
$ws.Range("K49").Value = 416
$ws.Range("M49").Value = "Yes"

# 3) Rename the required_quest_item_id of the "Alchemy Corrupted Church" row
#    (now at row 55) from "Purgatory's Lantern" to "Key to The Abandonded Church".
$ws.Range("D55").Value = "Key to The Abandonded Church"

# 4) Remove the now-duplicated old "Abandonded Chapel" row, which the insert
#    above pushed down to row 57.
$ws.Rows.Item(57).Delete()

# Re-apply best-fit sizing to column D to reflect the new (longer) text
# ("Key to The Abandonded Church" no longer fits the old best-fit width).
$ws.Columns.Item(4).ColumnWidth = 33.25
